$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$np = $s.NotesPage
$np.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = ""
